$d = $word.ActiveDocument

# --- 1. Replace the italic meta-description text (near the bottom) with the new image prompt text ---
# (Use explicit Start/End Range.Text assignment rather than Find.Execute's replace, which
#  would trigger Word's smart-quotes AutoCorrect and mangle the straight apostrophe in "game's".)
$oldDesc = "Experience Book of Stars online slot game for free with this review. Learn about the game design, mechanics, bonus features, and how it compares to similar slots."
$newDesc = "Prompt: Create a cartoon-style feature image for Book of Stars that features a happy Maya warrior with glasses. The image should showcase the Maya warrior standing in front of the game grid, with constellations and stars framing the top of it. The warrior should be holding the Book of Stars, the wild symbol in the game, with a big smile on his face as he looks out at the viewer. He should be wearing glasses to show his intelligence and add a touch of humor. The background should showcase the astral theme and include symbols with gold details, like playing card symbols, to represent the game's design features. The overall style should be colorful and eye-catching to attract potential players."

$descFound = $false
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    If ($para.Range.Text.TrimEnd([char]13, [char]7) -eq $oldDesc) {
        $r = $d.Range($para.Range.Start, $para.Range.End)
        $r.Text = $newDesc
        $descFound = $true
        Break
    }
}

# --- 2. Remove the old bold "Play Book of Stars..." paragraph near the bottom ---
$oldTitleText = "Play Book of Stars Free | Review of Online Slot Game"
$found = $false
For ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    If ($para.Range.Text.TrimEnd([char]13, [char]7) -eq $oldTitleText -and $para.Range.Bold) {
        $para.Range.Delete()
        $found = $true
        Break
    }
}

# --- 3. Insert a new "Meta description" paragraph right after the title (Heading1) ---
$titlePara = $d.Paragraphs(1)
$newParaRange = $titlePara.Range.InsertParagraphAfter()
$target = $d.Paragraphs(2).Range

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Experience Book of Stars online slot game for free with this review. Learn about the game design, mechanics, bonus features, and how it compares to similar slots.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($metaXml)

Write-Output "removed_old_title=$found description_replaced=$descFound"
